$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# USD Amount (column T) for the row-2 transaction was corrected.
$ws.Range("T2").Value = 204937
